# The deck originally shipped with its slide-master theme (theme1.xml,
# "Integral" / "Red Violet" colour scheme) and its notes-master theme
# (theme2.xml, the stock "Office Theme" colour scheme).
#
# The authored change swaps the two: the presentation's theme becomes the
# default Office colour scheme, while the old Integral/Red Violet palette
# is retired to the notes-master slot.
#
# Apply the new (Office Theme) palette to the presentation's theme by
# rewriting every slot of the 12-colour theme colour scheme -- this is
# the PowerPoint COM surface for editing a:clrScheme entries (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) on the live theme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RGB() equivalent values (R + G*256 + B*65536) for the Office Theme
# palette, in clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
